# Auto-generated edit script: refresh crypto price/volume table cells
# to match the updated data snapshot (GitHub Actions update).
#
# Numeric-looking values are written with a leading apostrophe so Excel
# stores them as literal text (preserving trailing zeros / multi-dot
# "thousands" formatting such as "26.872.19") instead of auto-converting
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.872.19'
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').Value = '1.809.84'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''309.89'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '''0.4645'
$ws.Range('E7').Value = '  +4.12%  '
$ws.Range('D8').Value = '''0.3705'
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('D9').Value = '''0.07362'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('D10').Value = '''0.8755'
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('D11').Value = '''20.46'
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('D12').Value = '1.846.07'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').Value = '''5.358'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '''6.507'
$ws.Range('E14').Value = '  -3.17%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = '''91.79'
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').Value = '''0.07044'
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = '''0.000008693'
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '''14.73'
$ws.Range('E20').Value = '  -2.36%  '
$ws.Range('D21').Value = '26.877.61'
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('D22').Value = '''5.312'
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('D23').Value = '''10.59'
$ws.Range('E23').Value = '  -3.16%  '
$ws.Range('D24').Value = '2.018.64'
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('D25').Value = '''1.896'
$ws.Range('E25').Value = '  -3.06%  '
$ws.Range('D26').Value = '''151.55'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').Value = '''18.39'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').Value = '''2.153'
$ws.Range('E28').Value = '  -6.19%  '
$ws.Range('D29').Value = '''5.330'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').Value = '''115.92'
$ws.Range('E30').Value = '  -0.99%  '
$ws.Range('D31').Value = '''0.08902'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = '''0.7540'
$ws.Range('E32').Value = '  -5.29%  '
$ws.Range('D33').Value = '''1.160'
$ws.Range('E33').Value = '  -3.26%  '
$ws.Range('D34').Value = '''4.458'
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('E35').Value = '  -1.56%  '
$ws.Range('D36').Value = '''1.001'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '''1.102'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('D38').Value = '''0.01968'
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('D39').Value = '''2.445'
$ws.Range('E39').Value = '  +4.58%  '
$ws.Range('D40').Value = '''0.05242'
$ws.Range('E40').Value = '  -0.62%  '
$ws.Range('D41').Value = '''2.925'
$ws.Range('E41').Value = '  +1.82%  '
$ws.Range('D42').Value = '''0.5315'
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('D43').Value = '''7.177'
$ws.Range('E43').Value = '  -1.66%  '
$ws.Range('D44').Value = '''0.1664'
$ws.Range('E44').Value = '  -2.36%  '
$ws.Range('D45').Value = '''8.490'
$ws.Range('E45').Value = '  -2.44%  '
$ws.Range('D46').Value = '''0.4979'
$ws.Range('E46').Value = '  -1.98%  '
$ws.Range('D47').Value = '''10.30'
$ws.Range('E47').Value = '  -3.31%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').Value = '''1.001'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''103.88'
$ws.Range('E49').Value = '  -1.50%  '
$ws.Range('D50').Value = '''1.668'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').Value = '''0.06297'
$ws.Range('E51').Value = '  -1.54%  '
